# Update "想去人数" (F column) values on 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 433
    3  = 1452
    4  = 949
    6  = 2119
    8  = 1282
    10 = 117
    11 = 36
    12 = 306
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
